$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generations 0-3 (rows 2-5) -> 7728
$ws.Range("C2:C5").Value = 7728

# Generations 4-250 (rows 6-252) -> 7310
$ws.Range("C6:C252").Value = 7310
